# Update the L5CG9_routine_3rdSem worksheet with new 2nd-semester schedule.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (A1:L1) - Room/Block/Group reordered, Level and Course added.
$headers = @("Day","Time","Module Code","Module Title","Hours","Class Type","Lecturer","Room","Block","Group","Level","Course")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# New data rows 2-10.
$data = @(
    @("SUN","13:00-15:00","5CS022","Human Computer Interaction",2,"Tutorial","Mr. Pravash Karki","SR-01 Bantok","WLV","L5CG9",5,"BCS"),
    @("MON","9:00-11:30","5CS022","Human Computer Interaction",2.5,"Workshop","Mr. Pravash Karki","SR-03 Wolves","WLV","L5CG9",5,"BCS"),
    @("TUE","9:30-11:30","5CS024","Collaborative Development",2,"Lecture","Mr. Udaya Kandel","LT-03 Walsall","WLV","L5CG(9+10+11)",5,"BCS"),
    @("TUE","12:00-14:00","5CS022","Human Computer Interaction",2,"Lecture","Mr. Pravash Karki","LT-02 Telford","WLV","L5CG(9+10+11)",5,"BCS"),
    @("WED","9:00-11:00","5CS020","Distributed and Cloud Systems Programming",2,"Lecture","Mr. Sumanta Silwal","LT-02 Telford","WLV","L5CG(9+10+11)",5,"BCS"),
    @("THU","9:30-11:30","5CS024","Collaborative Development",2,"Tutorial","Mr. Biraj Dulal","SR-02 Bilston","WLV","L5CG9",5,"BCS"),
    @("THU","13:00-15:00","5CS020","Distributed and Cloud Systems Programming",2,"Tutorial","Mr. Sumanta Silwal","TR-02 Stafford","WLV","L5CG9",5,"BCS"),
    @("FRI","9:30-12:00","5CS020","Distributed and Cloud Systems Programming",2.5,"Workshop","Mr. Sumanta Silwal","Lab-01 Mander","WLV","L5CG9",5,"BCS"),
    @("FRI","12:30-15:30","5CS024","Collaborative Development",2.5,"Workshop","Mr. Biraj Dulal","TR-05 Ranipokhari","HCK","L5CG9",5,"BCS")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
